# Applies the Behemoth_Profits commodity-pricing update across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR worksheets, matching the
# scheduled-runner recalculation of leve crafting costs/profits.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 17
$ws.Range("H17").Value = 1800
$ws.Range("J17").Value = 1800
$ws.Range("L17").Value = 5400
$ws.Range("N17").Value = -5736

# Row 74
$ws.Range("H74").Value = 3778
$ws.Range("I74").Value = 3857.4285
$ws.Range("K74").Value = 3857.4285
$ws.Range("M74").Value = -2921.4285

# Row 77
$ws.Range("H77").Value = 3778
$ws.Range("I77").Value = 3857.4285
$ws.Range("K77").Value = 19287.1425
$ws.Range("M77").Value = -14607.1425

# Row 98
$ws.Range("H98").Value = 52633588
$ws.Range("I98").Value = 71430240
$ws.Range("J98").Value = 2960
$ws.Range("K98").Value = 71430240
$ws.Range("L98").Value = 2960
$ws.Range("M98").Value = -71428742
$ws.Range("N98").Value = -5956

# Row 113
$ws.Range("H113").Value = 83335200
$ws.Range("J113").Value = 100002104
$ws.Range("L113").Value = 100002104
$ws.Range("N113").Value = -100008612

# Row 122
$ws.Range("H122").Value = 52633588
$ws.Range("I122").Value = 71430240
$ws.Range("J122").Value = 2960
$ws.Range("K122").Value = 214290720
$ws.Range("L122").Value = 8880
$ws.Range("M122").Value = -214288270
$ws.Range("N122").Value = -13780

# Row 127
$ws.Range("H127").Value = 733
$ws.Range("I127").Value = 733
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2199
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 2761
$ws.Range("N127").ClearContents()

# Row 129
$ws.Range("H129").Value = 1409.5555
$ws.Range("I129").Value = 1170.8572
$ws.Range("K129").Value = 3512.5716
$ws.Range("M129").Value = 1487.4284

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Row 135
$ws.Range("H135").Value = 1593.9286
$ws.Range("J135").Value = 2139.4
$ws.Range("L135").Value = 19254.6
$ws.Range("N135").Value = -24324.6

# Row 137
$ws.Range("H137").Value = 5393.3145
$ws.Range("J137").Value = 6120
$ws.Range("L137").Value = 18360
$ws.Range("N137").Value = -23460

$ws = $wb.Worksheets("ARM")
# Row 32
$ws.Range("H32").Value = 20840320
$ws.Range("I32").Value = 26318298
$ws.Range("K32").Value = 26318298
$ws.Range("M32").Value = -26318011

# Row 44
$ws.Range("H44").Value = 41285.855
$ws.Range("J44").Value = 41285.855
$ws.Range("L44").Value = 41285.855
$ws.Range("N44").Value = -42261.855

# Row 45
$ws.Range("H45").Value = 17244134
$ws.Range("I45").Value = 27779712
$ws.Range("J45").Value = 4098.364
$ws.Range("K45").Value = 27779712
$ws.Range("L45").Value = 4098.364
$ws.Range("M45").Value = -27779335
$ws.Range("N45").Value = -4852.364

# Row 61
$ws.Range("H61").Value = 37507050
$ws.Range("I61").Value = 62505660
$ws.Range("K61").Value = 62505660
$ws.Range("M61").Value = -62505448

# Row 63
$ws.Range("H63").Value = 6169.6924
$ws.Range("I63").Value = 2365.6667
$ws.Range("K63").Value = 2365.6667
$ws.Range("M63").Value = -1679.6667

# Row 66
$ws.Range("H66").Value = 6169.6924
$ws.Range("I66").Value = 2365.6667
$ws.Range("K66").Value = 11828.3335
$ws.Range("M66").Value = -8396.333500000001

# Row 74
$ws.Range("H74").Value = 15295769
$ws.Range("I74").Value = 20834862
$ws.Range("K74").Value = 20834862
$ws.Range("M74").Value = -20833988

# Row 77
$ws.Range("H77").Value = 15295769
$ws.Range("I77").Value = 20834862
$ws.Range("K77").Value = 104174310
$ws.Range("M77").Value = -104169942

# Row 136
$ws.Range("H136").Value = 37507050
$ws.Range("I136").Value = 62505660
$ws.Range("K136").Value = 187516980
$ws.Range("M136").Value = -187514430

$ws = $wb.Worksheets("BSM")
# Row 20
$ws.Range("H20").Value = 3162.25
$ws.Range("I20").Value = 3599.8667
$ws.Range("K20").Value = 3599.8667
$ws.Range("M20").Value = -3352.8667

# Row 86
$ws.Range("H86").Value = 2392.182
$ws.Range("I86").Value = 2994.7144
$ws.Range("J86").Value = 1337.75
$ws.Range("K86").Value = 2994.7144
$ws.Range("L86").Value = 1337.75
$ws.Range("M86").Value = -1871.7144
$ws.Range("N86").Value = -3583.75

# Row 89
$ws.Range("H89").Value = 2392.182
$ws.Range("I89").Value = 2994.7144
$ws.Range("J89").Value = 1337.75
$ws.Range("K89").Value = 14973.572
$ws.Range("L89").Value = 6688.75
$ws.Range("M89").Value = -9357.572
$ws.Range("N89").Value = -17920.75

# Row 134
$ws.Range("H134").Value = 4003468
$ws.Range("I134").Value = 1749.8572
$ws.Range("K134").Value = 5249.571599999999
$ws.Range("M134").Value = -2714.571599999999

$ws = $wb.Worksheets("CRP")
# Row 9
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10336

# Row 111
$ws.Range("H111").Value = 99989
$ws.Range("J111").Value = 99989
$ws.Range("L111").Value = 99989
$ws.Range("N111").Value = -108169

# Row 132
$ws.Range("H132").Value = 3101.0908
$ws.Range("I132").Value = 2436.158
$ws.Range("K132").Value = 7308.474
$ws.Range("M132").Value = -4778.474

# Row 134
$ws.Range("H134").Value = 2839.8276
$ws.Range("I134").Value = 2875.2307
$ws.Range("J134").Value = 2533
$ws.Range("K134").Value = 8625.6921
$ws.Range("L134").Value = 7599
$ws.Range("M134").Value = -6090.6921
$ws.Range("N134").Value = -12669

$ws = $wb.Worksheets("CUL")
# Row 109
$ws.Range("H109").Value = 1789.625
$ws.Range("J109").Value = 3500
$ws.Range("L109").Value = 10500
$ws.Range("N109").Value = -12580

$ws = $wb.Worksheets("GSM")
# Row 2
$ws.Range("H2").Value = 282.22223
$ws.Range("I2").Value = 56.666668
$ws.Range("J2").Value = 395
$ws.Range("K2").Value = 56.666668
$ws.Range("L2").Value = 395
$ws.Range("M2").Value = 56.333332
$ws.Range("N2").Value = -621

# Row 102
$ws.Range("H102").Value = 2726.5
$ws.Range("I102").Value = 2260.4666
$ws.Range("K102").Value = 2260.4666
$ws.Range("M102").Value = -638.4666000000002

# Row 122
$ws.Range("H122").Value = 2517.5715
$ws.Range("I122").Value = 2156
$ws.Range("K122").Value = 6468
$ws.Range("M122").Value = -4018

# Row 132
$ws.Range("H132").Value = 34490812
$ws.Range("I132").Value = 50004056
$ws.Range("K132").Value = 150012168
$ws.Range("M132").Value = -150009638

$ws = $wb.Worksheets("LTW")
# Row 61
$ws.Range("H61").Value = 2188.5557
$ws.Range("I61").Value = 2188.5557
$ws.Range("K61").Value = 2188.5557
$ws.Range("M61").Value = -1986.5557

# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 112
$ws.Range("H112").Value = 91459
$ws.Range("J112").Value = 91459
$ws.Range("L112").Value = 91459
$ws.Range("N112").Value = -94413

# Row 113
$ws.Range("H113").Value = 2188.5557
$ws.Range("I113").Value = 2188.5557
$ws.Range("K113").Value = 2188.5557
$ws.Range("M113").Value = -18.55569999999989

# Row 122
$ws.Range("H122").Value = 5929.923
$ws.Range("I122").Value = 5429
$ws.Range("K122").Value = 16287
$ws.Range("M122").Value = -13837

# Row 136
$ws.Range("H136").Value = 124539.38
$ws.Range("I136").Value = 26801.6
$ws.Range("J136").Value = 185625.5
$ws.Range("K136").Value = 80404.79999999999
$ws.Range("L136").Value = 556876.5
$ws.Range("M136").Value = -77854.79999999999
$ws.Range("N136").Value = -561976.5

$ws = $wb.Worksheets("WVR")
# Row 75
$ws.Range("H75").Value = 100119.664
$ws.Range("J75").Value = 100119.664
$ws.Range("L75").Value = 100119.664
$ws.Range("N75").Value = -101991.664

# Row 78
$ws.Range("H78").Value = 100119.664
$ws.Range("J78").Value = 100119.664
$ws.Range("L78").Value = 300358.992
$ws.Range("N78").Value = -309718.992

# Row 122
$ws.Range("H122").Value = 2485.394
$ws.Range("I122").Value = 2457.8572
$ws.Range("K122").Value = 7373.571599999999
$ws.Range("M122").Value = -4923.571599999999
